# sprint 5 all changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 2
# ------------------------------------------------------------------
# A2 used to be the bare number 753011515 with no special formatting.
# It becomes the TEXT value "755841007" with the same bordered style
# the rest of the data cells use. Copy that look from B3 first, then
# force text storage via NumberFormat "@" before writing the value so
# it isn't silently re-parsed back into a number.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Value = "755841007"

# F2: PREPAID -> PREPAID (OR) POSTPAID (format / style unchanged)
$ws.Range("F2").Value = "PREPAID (OR) POSTPAID"

# G2: brand-new blank cell, date-style number format with a thin border
$ws.Range("G2").Borders.LineStyle = 1
$ws.Range("G2").NumberFormat = "m/d/yyyy"

# ------------------------------------------------------------------
# Row 3
# ------------------------------------------------------------------
# A3: 755843100 -> 755841007, now picks up the same bordered style as
# the rest of row 3 (it had no style at all before). Write the numeric
# value BEFORE pasting the format so it stays a true number (pasting a
# "@"/Text-formatted style over an already-typed value does not turn it
# into text the way it would if the format existed first).
$ws.Range("A3").Value = 755841007
$ws.Range("B3").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null

# F3: POSTPAID -> PREPAID, now picks up the centred/bordered style used
# by the other cells in its row (it had no style at all before).
$ws.Range("F3").Value = "PREPAID"
$ws.Range("E3").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4122) | Out-Null

# G3: brand-new blank cell, same look as G2
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G3").PasteSpecial(-4122) | Out-Null

# ------------------------------------------------------------------
# Row 4 - brand-new row, mirrors row 3's shape/format exactly except
# for the MSISDN (A4) and Amount (E4) values.
# A4 is written BEFORE the format copy so it keeps its true numeric
# type; the other cells are written AFTER so they pick up row 3's
# Text ("@") number format the same way row 3 / row 2 store their
# "10"/"3"/"100" figures as text.
# ------------------------------------------------------------------
$ws.Range("A4").Value = 755841651

$ws.Range("A3:G3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null

$ws.Range("B4").Value = "Weekend Offer"
$ws.Range("C4").Value = "Unlimited Youtube for 7 days"
$ws.Range("D4").Value = "3"
$ws.Range("E4").Value = "99"
$ws.Range("F4").Value = "PREPAID"

# ------------------------------------------------------------------
# Sheet-level bits
# ------------------------------------------------------------------
$ws.Range("C2").Select() | Out-Null
